$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.749.35"
$ws.Range("E2").Value = "  +0.31%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.603.38"
$ws.Range("E3").Value = "  +0.42%  "

# Row 4
$ws.Range("E4").Value = "  +0.20%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "211.92"
$ws.Range("E5").Value = "  +0.20%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.514"
$ws.Range("E6").Value = "  +0.33%  "

# Row 7
$ws.Range("E7").Value = "  +0.20%  "

# Row 8
$ws.Range("E8").Value = "  +0.20%  "

# Row 9
$ws.Range("E9").Value = "  +0.26%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.62"
$ws.Range("E10").Value = "  +0.68%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0848"
$ws.Range("E11").Value = "  +0.86%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.828.10"
$ws.Range("E12").Value = "  +0.36%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.620.54"
$ws.Range("E13").Value = "  +2.02%  "

# Row 14
$ws.Range("E14").Value = "  +1.04%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.526"
$ws.Range("E15").Value = "  +0.54%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.13"
$ws.Range("E16").Value = "  +0.10%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0₃0740"
$ws.Range("E17").Value = "  -0.02%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "209.68"
$ws.Range("E18").Value = "  +0.25%  "

# Row 19
$ws.Range("E19").Value = "  +0.24%  "

# Row 20
$ws.Range("E20").Value = "  +1.52%  "

# Row 21
$ws.Range("E21").Value = "  +0.33%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.22"
$ws.Range("E22").Value = "  -5.01%  "

# Row 23
$ws.Range("E23").Value = "  +0.76%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "143.82"
$ws.Range("E24").Value = "  +0.48%  "

# Row 25
$ws.Range("E25").Value = "  +0.45%  "

# Row 26
$ws.Range("E26").Value = "  -0.32%  "

# Row 27
$ws.Range("E27").Value = "  -0.03%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.37"
$ws.Range("E28").Value = "  +0.29%  "

# Row 29
$ws.Range("E29").Value = "  -1.01%  "

# Row 30
$ws.Range("E30").Value = "  +0.27%  "

# Row 31
$ws.Range("E31").Value = "  +1.31%  "

# Row 32
$ws.Range("E32").Value = "  +0.77%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.287.58"
$ws.Range("E33").Value = "  -0.13%  "

# Row 34
$ws.Range("E34").Value = "  +1.15%  "

# Row 35
$ws.Range("E35").Value = "  +18.12%  "

# Row 36
$ws.Range("E36").Value = "  +0.41%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.586"
$ws.Range("E37").Value = "  -5.28%  "

# Row 38
$ws.Range("E38").Value = "  -0.36%  "

# Row 39
$ws.Range("E39").Value = "  -0.06%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.44"
$ws.Range("E40").Value = "  -0.25%  "

# Row 41
$ws.Range("E41").Value = "  +0.01%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.779"
$ws.Range("E42").Value = "  -0.32%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "62.67"
$ws.Range("E43").Value = "  -0.85%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.739.61"
$ws.Range("E44").Value = "  +0.43%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "90.44"
$ws.Range("E45").Value = "  -0.58%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.57"
$ws.Range("E46").Value = "  +0.46%  "

# Row 47
$ws.Range("E47").Value = "  +1.02%  "

# Row 48
$ws.Range("E48").Value = "  +0.76%  "

# Row 49
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.52"
$ws.Range("E49").Value = "  +2.06%  "

# Row 50
$ws.Range("B50").Value = "USDD"
$ws.Range("C50").Value = "https://coinranking.com/coin/z2PZIKQL7+usdd-usdd"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.00"
$ws.Range("E50").Value = "  +0.07%  "

# Row 51
$ws.Range("B51").Value = "Mantle"
$ws.Range("C51").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.399"
$ws.Range("E51").Value = "  +1.88%  "
